# fixed #113 OneWayRPG-113 (defense-up) - add the "food insufficient"
# message row that was missing from the "message" sheet.
#
# Mirrors the existing "out of money" row (row 75, id 73) with a new
# row 76 (id 74) right below it, using the same "yellow" message color.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("message")

$srcRow = 75
$newRow = 76

# --- values / formula -------------------------------------------------
$ws.Cells.Item($newRow, 1).Formula = "=ROW()-2"
$ws.Cells.Item($newRow, 2).Value   = "食糧が足りない"
$ws.Cells.Item($newRow, 3).Value   = "yellow"

# --- formatting: mirror the row above (same table style) --------------
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item($srcRow).RowHeight

for ($col = 1; $col -le 3; $col++) {
    $src = $ws.Cells.Item($srcRow, $col)
    $dst = $ws.Cells.Item($newRow, $col)

    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.NumberFormat = $src.NumberFormat

    if ($src.Interior.ColorIndex -ne -4142) {
        $dst.Interior.Pattern = $src.Interior.Pattern
        $dst.Interior.ColorIndex = $src.Interior.ColorIndex
    }

    foreach ($edge in 7, 8, 9, 10) {
        $dst.Borders.Item($edge).LineStyle = 1
        $dst.Borders.Item($edge).Weight = 2
    }
}
